# Append 9 additional rows of test data to the
# "master-reg_center_user_machine_" sheet (rows 22-30), matching the
# existing pattern used in rows 2-21, then restore the cursor/selection
# and page orientation as left by the author after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id, usr_id, machine_id triples for the new rows
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]          # regcntr_id
    $ws.Cells.Item($r, 2).Value = $row[1]          # usr_id
    $ws.Cells.Item($r, 3).Value = $row[2]          # machine_id
    $ws.Cells.Item($r, 4).Value = "eng"            # lang_code
    $ws.Cells.Item($r, 5).Value = $true            # is_active
    $ws.Cells.Item($r, 6).Value = "superadmin"     # cr_by
    $ws.Cells.Item($r, 7).Value = "now()"          # cr_dtimes
    $ws.Cells.Item($r, 8).Value = "now()"          # eff_dtimes
    $r++
}

# Set page orientation to portrait (as reflected in the saved pageSetup)
$ws.PageSetup.Orientation = 1

# Leave the active cell/selection where the author left it
$ws.Range("F14").Select()

$wb.Save()
